# Leaderboard update: add three new rows (10-12) with customer/salesperson
# data, extending the table that previously had blank placeholder rows
# 10 and 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 already carries the correct formatting (text style for A:C,
# date style for D) for the new rows, so clone it into row 12 before
# filling in values. This keeps the existing cell styles (and avoids
# Excel re-guessing number formats / types for a brand new row).
$ws.Range("A11:D11").Copy($ws.Range("A12:D12"))

# Row 10
$ws.Range("A10").Value = "Isaac's Pub"
$ws.Range("B10").Value = "Frisch, Isaac"
$ws.Range("C10").Value = "013"
$ws.Range("D10").Value = 45848

# Row 11
$ws.Range("A11").Value = "State Street"
$ws.Range("B11").Value = "Frisch, Isaac"
$ws.Range("C11").Value = "013"
$ws.Range("D11").Value = 45848

# Row 12 (new)
$ws.Range("A12").Value = "Reciprical Show"
$ws.Range("B12").Value = "Frisch, Isaac"
$ws.Range("C12").Value = "013"
$ws.Range("D12").Value = 45848

$ws.Range("B12").Select()
